$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.646.53"
Set-TextValue $ws.Range("E2") "  -1.78%  "
Set-TextValue $ws.Range("D3") "1.786.59"
Set-TextValue $ws.Range("E3") "  -1.82%  "
Set-TextValue $ws.Range("E4") "  +0.29%  "
Set-TextValue $ws.Range("D5") "308.09"
Set-TextValue $ws.Range("E5") "  -1.28%  "
Set-TextValue $ws.Range("E6") "  +0.33%  "
Set-TextValue $ws.Range("E7") "  +1.78%  "
Set-TextValue $ws.Range("D8") "0.3688"
Set-TextValue $ws.Range("E8") "  -1.47%  "
Set-TextValue $ws.Range("D9") "0.07241"
Set-TextValue $ws.Range("E9") "  -2.87%  "
Set-TextValue $ws.Range("D10") "0.8490"
Set-TextValue $ws.Range("E10") "  -2.58%  "
Set-TextValue $ws.Range("D11") "20.34"
Set-TextValue $ws.Range("E11") "  -2.77%  "
Set-TextValue $ws.Range("D12") "1.798.04"
Set-TextValue $ws.Range("E12") "  -1.08%  "
Set-TextValue $ws.Range("D13") "6.487"
Set-TextValue $ws.Range("E13") "  -3.91%  "
Set-TextValue $ws.Range("D14") "5.276"
Set-TextValue $ws.Range("E14") "  -1.21%  "
Set-TextValue $ws.Range("D15") "0.07015"
Set-TextValue $ws.Range("E15") "  -1.03%  "
Set-TextValue $ws.Range("D16") "90.42"
Set-TextValue $ws.Range("E16") "  -4.30%  "
Set-TextValue $ws.Range("E17") "  +0.37%  "
Set-TextValue $ws.Range("D18") "0.000008565"
Set-TextValue $ws.Range("E18") "  -2.12%  "
Set-TextValue $ws.Range("E19") "  +0.25%  "
Set-TextValue $ws.Range("D20") "14.57"
Set-TextValue $ws.Range("E20") "  -2.91%  "
Set-TextValue $ws.Range("D21") "26.649.97"
Set-TextValue $ws.Range("E21") "  -1.82%  "
Set-TextValue $ws.Range("D22") "5.244"
Set-TextValue $ws.Range("E22") "  +0.43%  "
Set-TextValue $ws.Range("D23") "10.61"
Set-TextValue $ws.Range("E23") "  -3.07%  "
Set-TextValue $ws.Range("D24") "2.012.83"
Set-TextValue $ws.Range("E24") "  -1.47%  "
Set-TextValue $ws.Range("E25") "  -4.20%  "
Set-TextValue $ws.Range("D26") "149.72"
Set-TextValue $ws.Range("E26") "  -1.22%  "
Set-TextValue $ws.Range("D27") "2.144"
Set-TextValue $ws.Range("E27") "  -11.63%  "
Set-TextValue $ws.Range("D28") "18.11"
Set-TextValue $ws.Range("E28") "  -2.41%  "
Set-TextValue $ws.Range("D29") "5.171"
Set-TextValue $ws.Range("E29") "  -2.90%  "
Set-TextValue $ws.Range("D30") "113.77"
Set-TextValue $ws.Range("E30") "  -3.62%  "
Set-TextValue $ws.Range("D31") "0.08824"
Set-TextValue $ws.Range("E31") "  +0.12%  "
Set-TextValue $ws.Range("D32") "0.7505"
Set-TextValue $ws.Range("E32") "  -2.08%  "
Set-TextValue $ws.Range("D33") "1.153"
Set-TextValue $ws.Range("E33") "  -1.96%  "
Set-TextValue $ws.Range("D34") "4.422"
Set-TextValue $ws.Range("E34") "  -3.09%  "
Set-TextValue $ws.Range("D35") "2.867"
Set-TextValue $ws.Range("E35") "  -0.60%  "
Set-TextValue $ws.Range("E36") "  +0.29%  "
Set-TextValue $ws.Range("D37") "1.110"
Set-TextValue $ws.Range("E37") "  +0.82%  "
Set-TextValue $ws.Range("E38") "  -2.57%  "
Set-TextValue $ws.Range("D39") "0.05193"
Set-TextValue $ws.Range("E39") "  -1.65%  "
Set-TextValue $ws.Range("D40") "7.110"
Set-TextValue $ws.Range("E40") "  -4.29%  "
Set-TextValue $ws.Range("E41") "  +0.16%  "
Set-TextValue $ws.Range("E42") "  +7.08%  "
Set-TextValue $ws.Range("D43") "0.5180"
Set-TextValue $ws.Range("E43") "  -2.50%  "
Set-TextValue $ws.Range("E44") "  -4.32%  "
Set-TextValue $ws.Range("D45") "8.427"
Set-TextValue $ws.Range("E45") "  -3.68%  "
Set-TextValue $ws.Range("D46") "0.4919"
Set-TextValue $ws.Range("E46") "  -2.77%  "
Set-TextValue $ws.Range("D47") "10.22"
Set-TextValue $ws.Range("E47") "  -3.42%  "
Set-TextValue $ws.Range("D48") "1.002"
Set-TextValue $ws.Range("E48") "  +0.41%  "
Set-TextValue $ws.Range("D49") "103.35"
Set-TextValue $ws.Range("E49") "  -2.19%  "
Set-TextValue $ws.Range("D50") "1.637"
Set-TextValue $ws.Range("E50") "  -4.08%  "
Set-TextValue $ws.Range("D51") "0.06272"
Set-TextValue $ws.Range("E51") "  -1.50%  "
